# Apply the betexplorer "Atualizado por script" update:
#  - A bunch of existing rows (4-9, 17, 21, 24-25, 27-31, 44-49, 53-55, 58,
#    61, 95-98, 105-107) get their match data (columns F..V) replaced with
#    the match data that used to live in a different row (a pure row-content
#    permutation - columns A..E, the index/date metadata, stay put).
#  - One brand-new match row (110) is appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (source row's F:V content moves into target row)
$mapping = @{
  4=8; 5=7; 6=9; 7=5; 8=4; 9=6;
  17=21; 21=17;
  24=25; 25=24;
  27=30; 28=27; 29=31; 30=28; 31=29;
  44=49; 45=44; 46=45; 47=46; 48=47; 49=48;
  53=54; 54=55; 55=53;
  58=61; 61=58;
  95=97; 96=98; 97=95; 98=96;
  105=106; 106=107; 107=105
}

$cols = 6..22  # F (home) .. V (betexplorer url)

# Phase 1: snapshot the CURRENT (pre-edit) F:V values of every row that is
# involved in the permutation, before any writes happen (several rows both
# give data to, and receive data from, other rows in the same cycle).
$snapshot = @{}
foreach ($r in $mapping.Keys) {
  $rowVals = @{}
  foreach ($c in $cols) {
    $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
  }
  $snapshot[$r] = $rowVals
}

# Phase 2: write each target row's F:V cells from the snapshot of its
# mapped source row.
foreach ($r in $mapping.Keys) {
  $src = $mapping[$r]
  $rowVals = $snapshot[$src]
  foreach ($c in $cols) {
    $ws.Cells.Item($r, $c).Value = $rowVals[$c]
  }
}

# Append the new match row (110), matching formatting of the preceding
# data row (109) - bold/bordered index cell, date-formatted date cell.
$ws.Range("A109:V109").Copy($ws.Range("A110:V110"))

$ws.Cells.Item(110, 1).Value = 109
$ws.Cells.Item(110, 2).Value = "portugal"
$ws.Cells.Item(110, 3).Value = "liga-3"
$ws.Cells.Item(110, 4).Value = "2023-2024"
$ws.Cells.Item(110, 5).Value = 45248.66666666666
$ws.Cells.Item(110, 6).Value = "SC Vianense"
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = "AD Fafe"
$ws.Cells.Item(110, 9).Value = 1
$ws.Cells.Item(110, 10).Value = 2.09
$ws.Cells.Item(110, 11).Value = "14/11/2023 22:42"
$ws.Cells.Item(110, 12).Value = 2.44
$ws.Cells.Item(110, 13).Value = "18/11/2023 15:56"
$ws.Cells.Item(110, 14).Value = 3.33
$ws.Cells.Item(110, 15).Value = "14/11/2023 22:42"
$ws.Cells.Item(110, 16).Value = 3.18
$ws.Cells.Item(110, 17).Value = "18/11/2023 15:56"
$ws.Cells.Item(110, 18).Value = 3.41
$ws.Cells.Item(110, 19).Value = "14/11/2023 22:42"
$ws.Cells.Item(110, 20).Value = 3.12
$ws.Cells.Item(110, 21).Value = "18/11/2023 15:56"
$ws.Cells.Item(110, 22).Value = "https://www.betexplorer.com/football/portugal/liga-3/sc-vianense-ad-fafe/zLFbmQwc/"

Write-Host "Edit applied."
